$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells H1:J1 (matching result of OCR/master-data matching) ---
$ws.Range("H1").Value = "Coverage (raw)"
$ws.Range("I1").Value = "MatchScore"
$ws.Range("J1").Value = "Matched"

# Give the new headers the same look as the existing header row (bold,
# centered, bordered) by copying the formatting from A1.
$ws.Range("A1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

# --- Update the policy number on row 2 to the matched/suffixed value ---
$ws.Range("A2").Value = "IHA1801Y_13"

# --- New data cells H2:J2 ---
$ws.Range("H2").Value = "วัคซีนไข้หวัดใหญ่"

# Force these two to be stored as plain text (not auto-coerced to a
# number / boolean) using the leading-apostrophe text prefix, then reset
# the cell style to the default (unstyled) data-row look.
$ws.Range("I2").Value = "'100"
$ws.Range("I2").Style = $ws.Range("B2").Style

$ws.Range("J2").Value = "'TRUE"
$ws.Range("J2").Style = $ws.Range("B2").Style

$excel.CutCopyMode = $false
